$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header labels for each skill category
$ws.Range("A1").Value = "Core Foundations"
$ws.Range("H1").Value = "Backend"
$ws.Range("O1").Value = "Frontend"
$ws.Range("AA1").Value = "Databases & ORMs"
$ws.Range("AG1").Value = "DevOps & Cloud"

# Merge each category header across its column span
$ws.Range("A1:G1").Merge()
$ws.Range("H1:N1").Merge()
$ws.Range("O1:Z1").Merge()
$ws.Range("AA1:AF1").Merge()
$ws.Range("AG1:AK1").Merge()

# Format the whole header row: bold, larger font, centered, wrapped
$header = $ws.Range("A1:AK1")
$header.Font.Bold = $true
$header.Font.Size = 18
$header.HorizontalAlignment = -4108
$header.VerticalAlignment = -4108
$header.WrapText = $true

# Taller header row to fit wrapped text
$ws.Rows.Item(1).RowHeight = 42.6

# Print as portrait
$ws.PageSetup.Orientation = 1
